$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-26 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-27 Saturday", 2)

$d.Content.Find.Execute("306×9=", $true, $false, $false, $false, $false, $true, 1, $false, "752×7=", 2)
$d.Content.Find.Execute("516×5=", $true, $false, $false, $false, $false, $true, 1, $false, "988×8=", 2)
$d.Content.Find.Execute("142×2=", $true, $false, $false, $false, $false, $true, 1, $false, "444×8=", 2)
$d.Content.Find.Execute("601×5=", $true, $false, $false, $false, $false, $true, 1, $false, "892×5=", 2)
$d.Content.Find.Execute("835×2=", $true, $false, $false, $false, $false, $true, 1, $false, "998×4=", 2)

$d.Content.Find.Execute("529×8=", $true, $false, $false, $false, $false, $true, 1, $false, "389×6=", 2)
$d.Content.Find.Execute("305×6=", $true, $false, $false, $false, $false, $true, 1, $false, "278×7=", 2)
$d.Content.Find.Execute("235×4=", $true, $false, $false, $false, $false, $true, 1, $false, "943×8=", 2)
$d.Content.Find.Execute("620×7=", $true, $false, $false, $false, $false, $true, 1, $false, "494×2=", 2)
$d.Content.Find.Execute("238×9=", $true, $false, $false, $false, $false, $true, 1, $false, "432×9=", 2)

$d.Content.Find.Execute("838×3=", $true, $false, $false, $false, $false, $true, 1, $false, "260×2=", 2)
$d.Content.Find.Execute("834×8=", $true, $false, $false, $false, $false, $true, 1, $false, "788×2=", 2)
$d.Content.Find.Execute("723×9=", $true, $false, $false, $false, $false, $true, 1, $false, "466×6=", 2)
$d.Content.Find.Execute("845×3=", $true, $false, $false, $false, $false, $true, 1, $false, "643×8=", 2)
$d.Content.Find.Execute("338×8=", $true, $false, $false, $false, $false, $true, 1, $false, "369×4=", 2)

$d.Content.Find.Execute("231×3=", $true, $false, $false, $false, $false, $true, 1, $false, "866×6=", 2)
$d.Content.Find.Execute("279×3=", $true, $false, $false, $false, $false, $true, 1, $false, "895×8=", 2)
$d.Content.Find.Execute("906×5=", $true, $false, $false, $false, $false, $true, 1, $false, "983×5=", 2)
$d.Content.Find.Execute("483×6=", $true, $false, $false, $false, $false, $true, 1, $false, "335×6=", 2)
$d.Content.Find.Execute("878×6=", $true, $false, $false, $false, $false, $true, 1, $false, "786×4=", 2)

$d.Content.Find.Execute("849×7=", $true, $false, $false, $false, $false, $true, 1, $false, "878×9=", 2)
$d.Content.Find.Execute("868×7=", $true, $false, $false, $false, $false, $true, 1, $false, "582×8=", 2)
$d.Content.Find.Execute("507×4=", $true, $false, $false, $false, $false, $true, 1, $false, "480×2=", 2)
$d.Content.Find.Execute("574×3=", $true, $false, $false, $false, $false, $true, 1, $false, "246×8=", 2)
$d.Content.Find.Execute("955×8=", $true, $false, $false, $false, $false, $true, 1, $false, "274×9=", 2)
